$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 entirely (shifts nothing below it, just removes the row)
$ws.Rows.Item(3).Delete()

# Update row 2 values
$ws.Range("A2").Value = "một nửa"
$ws.Range("B2").Value = "a half"

# Apply a red fill to A2
$ws.Range("A2").Interior.Color = 255

# Update the active selection
$ws.Range("M5").Select()
